# Update Shikhar Dhawan's Delhi Capitals innings activity figures
# (runs / balls / fours / sixes columns C:F) to the latest values pulled
# into the Excel form. Only the cells whose figures actually changed are
# touched; everything else is left exactly as-is.
#
# The sheet stores these numbers as text (e.g. "54" rather than 54), so
# each new value is written as a string to avoid Excel silently
# re-typing it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column, new value (row/col are 1-based; C=3, D=4, E=5, F=6)
$updates = @(
    @(2,3,"54"),
    @(2,4,"41"),
    @(2,5,"6"),
    @(2,6,"0"),
    @(3,3,"0"),
    @(3,4,"2"),
    @(3,5,"0"),
    @(3,6,"0"),
    @(4,3,"0"),
    @(4,4,"1"),
    @(4,5,"0"),
    @(6,3,"0"),
    @(6,4,"2"),
    @(6,5,"0"),
    @(6,6,"0"),
    @(8,3,"57"),
    @(8,4,"33"),
    @(8,5,"6"),
    @(8,6,"2"),
    @(9,3,"101"),
    @(9,4,"58"),
    @(9,5,"14"),
    @(9,6,"1"),
    @(10,3,"69"),
    @(10,4,"52"),
    @(10,5,"6"),
    @(10,6,"1"),
    @(11,3,"106"),
    @(11,4,"61"),
    @(11,5,"12"),
    @(11,6,"3"),
    @(12,3,"15"),
    @(12,4,"13"),
    @(12,5,"3"),
    @(13,3,"34"),
    @(13,4,"31"),
    @(13,5,"4"),
    @(14,3,"26"),
    @(14,4,"16"),
    @(14,5,"2"),
    @(14,6,"2"),
    @(15,3,"35"),
    @(15,4,"27"),
    @(15,6,"1"),
    @(16,3,"0"),
    @(16,4,"2"),
    @(16,5,"0"),
    @(17,3,"32"),
    @(17,4,"28"),
    @(17,5,"3"),
    @(18,3,"5"),
    @(18,4,"4"),
    @(18,5,"1"),
    @(18,6,"0")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newValue = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    # Prefix with an apostrophe so Excel keeps the entry as text (matching
    # the existing text-stored numbers) instead of converting it to a number.
    $cell.Formula = "'" + $newValue
}

$wb.Save()
